$d = $word.ActiveDocument

# Locate the paragraph that currently starts with "Sohan" (the one that will be
# split into a new "FRNSW Visit to Viva Rosehill 5/2/20" paragraph followed by
# the original "Sohan - learned lots..." paragraph).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Sohan*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Sohan' paragraph"
}

$targetPara = $d.Paragraphs($targetIndex)

# Replace the whole paragraph (text + trailing bookmark) with two paragraphs:
#   1) "FRNSW Visit to Viva Rosehill 5/2/20" - carries the _GoBack bookmark
#   2) "Sohan" / " – learned lots about risk..." - same two runs as before
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>FRNSW Visit to Viva Rosehill 5/2/20</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>Sohan</w:t></w:r><w:r><w:t xml:space="preserve"> – learned lots about risk and what it is (product of frequency and consequence)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$paraCountBefore = $d.Paragraphs.Count
$targetPara.Range.InsertXML($xmlFrag)

# InsertXML leaves behind the original paragraph mark as a trailing empty
# paragraph once the 1-paragraph range has been replaced by 2 paragraphs;
# remove that leftover empty paragraph so the structure matches exactly
# (the "Sohan" run content followed directly by the section properties).
$paraCountAfter = $d.Paragraphs.Count
if ($paraCountAfter -gt ($paraCountBefore + 1)) {
    $sohanPara = $d.Paragraphs($targetIndex + 1)
    $trailingPara = $d.Paragraphs($targetIndex + 2)
    $cleanupRange = $d.Range($sohanPara.Range.End - 1, $trailingPara.Range.End)
    $cleanupRange.Delete()
}

Write-Output "Paragraph count:"
Write-Output $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i : $($d.Paragraphs($i).Range.Text)"
}
